$d = $word.ActiveDocument

# 1. Update the date text: 1-6-2014 -> 1-6-2015
$d.Content.Find.Execute("1-6-2014", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1-6-2015", 2)

# 2. Move the "_GoBack" bookmark so that it sits right after the date
#    run (collapsed, still inside that paragraph) instead of inside the
#    "Doel" paragraph's text, where it originally lived.

# Locate the paragraph that now holds the date.
$p = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*1-6-2015*") {
        $p = $para
        break
    }
}

$pEnd = $p.Range.End

# Insert a one-character placeholder immediately before the paragraph
# mark so we have a real (non-empty) range to anchor the bookmark to -
# a truly zero-length range built right at (paragraph end - 1) gets
# mis-placed by the engine.
$placeholder = $d.Range($pEnd - 1, $pEnd - 1)
$placeholder.InsertAfter("X")

# Wrap the bookmark tightly around that placeholder character. Adding a
# bookmark named "_GoBack" replaces/removes any existing bookmark with
# that same name elsewhere in the document (i.e. the old one in the
# "Doel" paragraph disappears automatically).
$pEnd2 = $p.Range.End
$markerRange = $d.Range($pEnd2 - 2, $pEnd2 - 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

# Delete the placeholder again; the bookmark collapses to a zero-length
# bookmark left sitting right after the date run, still inside the
# paragraph, matching the target layout.
$markerRange2 = $d.Range($pEnd2 - 2, $pEnd2 - 1)
$markerRange2.Text = ""
